# Update the cryptos price/volume table with the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values that often look like plain numbers
# (e.g. "1.00", "380.13"). Mark the whole column as Text up front so the
# new values are stored verbatim as strings, same as the scraper output,
# instead of being auto-coerced to numbers by Excel's input parser.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "51.099.39"
$ws.Range("E2").Value = "  +0.19%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.961.69"
$ws.Range("E3").Value = "  +0.81%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "380.13"
$ws.Range("E5").Value = "  +1.53%  "

# Row 6 - Solana
$ws.Range("D6").Value = "102.36"
$ws.Range("E6").Value = "  +0.51%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.544"
$ws.Range("E7").Value = "  +1.74%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.29%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "36.58"
$ws.Range("E10").Value = "  +0.17%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.04%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  +2.00%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.423.91"
$ws.Range("E13").Value = "  +0.87%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +5.96%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "18.33"
$ws.Range("E15").Value = "  +2.15%  "

# Row 16 - Uniswap
$ws.Range("D16").Value = "12.01"
$ws.Range("E16").Value = "  +67.45%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.970.12"
$ws.Range("E17").Value = "  +1.14%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +2.23%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "51.181.31"
$ws.Range("E19").Value = "  +0.43%  "

# Row 20 - ImmutableX
$ws.Range("E20").Value = "  -1.14%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").Value = "12.41"
$ws.Range("E21").Value = "  -1.18%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "70.14"
$ws.Range("E23").Value = "  +2.65%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "3.28"
$ws.Range("E24").Value = "  +13.84%  "

# Row 25 - BitcoinCash
$ws.Range("D25").Value = "267.97"
$ws.Range("E25").Value = "  +1.28%  "

# Row 26 - Filecoin
$ws.Range("D26").Value = "7.92"
$ws.Range("E26").Value = "  -2.65%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "7.23"
$ws.Range("E27").Value = "  -8.05%  "

# Row 28 - Dai
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.03%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  -0.49%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "25.90"
$ws.Range("E30").Value = "  +1.07%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.109"
$ws.Range("E31").Value = "  -2.15%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "10.44"
$ws.Range("E32").Value = "  +5.95%  "

# Row 33 - InjectiveProtocol
$ws.Range("D33").Value = "34.49"
$ws.Range("E33").Value = "  +2.61%  "

# Row 34/35 - Toncoin and OKB swap ranking positions
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "51.02"
$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").Value = "2.07"
$ws.Range("E35").Value = "  +2.71%  "

# Row 36 - VeChain
$ws.Range("E36").Value = "  -3.10%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.02%  "

# Row 38 - LidoDAOToken
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").Value = "  +9.10%  "

# Row 39 - Stellar
$ws.Range("E39").Value = "  +1.78%  "

# Row 40/41 - ARBITRUM and Celestia swap ranking positions
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "16.60"
$ws.Range("E40").Value = "  +1.37%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "1.83"
$ws.Range("E41").Value = "  +2.56%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -1.81%  "

# Row 43 - Monero
$ws.Range("D43").Value = "124.80"
$ws.Range("E43").Value = "  +3.52%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "21.65"
$ws.Range("E44").Value = "  +3.29%  "

# Row 45 - NEARProtocol
$ws.Range("D45").Value = "3.53"
$ws.Range("E45").Value = "  +9.26%  "

# Row 46/47 - WEMIXToken and ApeXProtocol swap ranking positions
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "2.38"
$ws.Range("E46").Value = "  +2.87%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "2.02"
$ws.Range("E47").Value = "  -1.28%  "

# Row 48/49 - TheGraph and Maker swap ranking positions
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.051.19"
$ws.Range("E48").Value = "  +3.94%  "

$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "0.269"
$ws.Range("E49").Value = "  -6.49%  "

# Row 50 - BEAM
$ws.Range("E50").Value = "  -6.89%  "

# Row 51 - THORChain
$ws.Range("E51").Value = "  +7.47%  "
